# Add "max_online" / "max_offline" columns to the processes sheet, matching
# the commit "Added option to define maximum online and offline times in
# the input data."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("processes")

# Insert two new columns before the old column M (initial_state), pushing
# initial_state/delay from M/N to O/P. Use xlFormatFromLeftOrAbove (-4161)
# so the inserted cells inherit column L's existing number format/style
# instead of creating a brand-new style entry.
$ws.Range("M1:N1").EntireColumn.Insert(-4161)

$ws.Range("M1").Value = "max_online"
$ws.Range("N1").Value = "max_offline"

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 13).Value = 0
    $ws.Cells.Item($r, 14).Value = 0
}

# Update the saved selection / active-sheet state to match the authored
# workbook: the "processes" sheet becomes the active tab (selection M12),
# and "process_topology" loses tabSelected (selection H22).
$ws5 = $wb.Worksheets.Item("process_topology")
$ws5.Range("H22").Select()

$ws.Range("M12").Select()
$ws.Activate()
